# Updated: st 14. 04. 2021
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revised AgTests (F) and AgPosit (G) figures for previously reported days
$ws.Range("F362").Value = 228917
$ws.Range("G362").Value = 3180
$ws.Range("F363").Value = 188507
$ws.Range("G363").Value = 2765
$ws.Range("F364").Value = 167899
$ws.Range("G364").Value = 2471
$ws.Range("F365").Value = 184025
$ws.Range("G365").Value = 2394
$ws.Range("F366").Value = 339944
$ws.Range("G366").Value = 2853
$ws.Range("F367").Value = 765955
$ws.Range("G367").Value = 3920
$ws.Range("G368").Value = 2271
$ws.Range("F369").Value = 233875
$ws.Range("F370").Value = 181954
$ws.Range("G370").Value = 2045
$ws.Range("F371").Value = 159763
$ws.Range("G371").Value = 1957
$ws.Range("F372").Value = 179602
$ws.Range("G372").Value = 1868
$ws.Range("F373").Value = 348854
$ws.Range("G373").Value = 2369
$ws.Range("F374").Value = 772266
$ws.Range("F375").Value = 348699
$ws.Range("G375").Value = 1848
$ws.Range("G376").Value = 2223
$ws.Range("F377").Value = 176715
$ws.Range("G377").Value = 1814
$ws.Range("F378").Value = 157194
$ws.Range("G378").Value = 1544
$ws.Range("F379").Value = 179275
$ws.Range("G379").Value = 1610
$ws.Range("F380").Value = 344195
$ws.Range("G380").Value = 2014
$ws.Range("G381").Value = 2683
$ws.Range("F382").Value = 357577
$ws.Range("G382").Value = 1573
$ws.Range("F383").Value = 220752
$ws.Range("G383").Value = 1766
$ws.Range("F384").Value = 171908
$ws.Range("G384").Value = 1495
$ws.Range("G385").Value = 1407
$ws.Range("F386").Value = 182205
$ws.Range("G386").Value = 1358
$ws.Range("F387").Value = 351630
$ws.Range("F388").Value = 728541
$ws.Range("G388").Value = 2196
$ws.Range("F389").Value = 353062
$ws.Range("G389").Value = 1303
$ws.Range("F390").Value = 220335
$ws.Range("G390").Value = 1516
$ws.Range("F391").Value = 176656
$ws.Range("G391").Value = 1191
$ws.Range("F392").Value = 220075
$ws.Range("G392").Value = 1216
$ws.Range("F393").Value = 298618
$ws.Range("G393").Value = 1194
$ws.Range("F395").Value = 738199
$ws.Range("G395").Value = 1918
$ws.Range("F398").Value = 291065
$ws.Range("G398").Value = 1440
$ws.Range("F399").Value = 194969
$ws.Range("G399").Value = 949
$ws.Range("F400").Value = 145089
$ws.Range("F401").Value = 262452
$ws.Range("F402").Value = 693162
$ws.Range("F403").Value = 333347
$ws.Range("G403").Value = 713
$ws.Range("F404").Value = 216072
$ws.Range("G404").Value = 884

# Append the new day's data as row 405
$ws.Range("A405").Value = 44299
$ws.Range("B405").Value = 373107
$ws.Range("C405").Value = 8593
$ws.Range("D405").Value = 1069
$ws.Range("E405").Value = 10798
$ws.Range("F405").Value = 145553
$ws.Range("G405").Value = 627
